$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1331
$ws.Range("F6").Value = 7660
$ws.Range("F7").Value = 95
$ws.Range("F9").Value = 2091
$ws.Range("F10").Value = 8451
$ws.Range("F12").Value = 52
$ws.Range("F14").Value = 5660
$ws.Range("F16").Value = 2619
$ws.Range("F17").Value = 1138
$ws.Range("F18").Value = 4597
$ws.Range("F20").Value = 403
$ws.Range("F21").Value = 94
$ws.Range("F23").Value = 532
$ws.Range("F24").Value = 3522
$ws.Range("F29").Value = 3020
$ws.Range("F30").Value = 37
$ws.Range("F31").Value = 105
$ws.Range("F32").Value = 344
$ws.Range("F34").Value = 311
$ws.Range("F35").Value = 402
$ws.Range("F36").Value = 656
$ws.Range("F38").Value = 883
$ws.Range("F39").Value = 1786
$ws.Range("F42").Value = 19
$ws.Range("F43").Value = 2892
$ws.Range("F49").Value = 1

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 106
$ws.Range("F3").Value = 124
$ws.Range("F5").Value = 50
$ws.Range("F6").Value = 4

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1325

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1325
$ws.Range("F5").Value = 1331
$ws.Range("F6").Value = 7660
$ws.Range("F7").Value = 95
$ws.Range("F9").Value = 2091
$ws.Range("F10").Value = 8451
$ws.Range("F12").Value = 52
$ws.Range("F14").Value = 5660
$ws.Range("F16").Value = 2619
$ws.Range("F17").Value = 1138
$ws.Range("F18").Value = 4597
$ws.Range("F19").Value = 403
$ws.Range("F20").Value = 106
$ws.Range("F22").Value = 124
$ws.Range("F23").Value = 532
$ws.Range("F25").Value = 3522
$ws.Range("F29").Value = 3020
$ws.Range("F30").Value = 344
$ws.Range("F32").Value = 311
$ws.Range("F33").Value = 50
$ws.Range("F34").Value = 402
$ws.Range("F35").Value = 656
$ws.Range("F37").Value = 883
$ws.Range("F39").Value = 1786
$ws.Range("F42").Value = 19
$ws.Range("F43").Value = 2892
